# Adds two new columns, I ("I0") and J ("IF"), to the sheet, mirroring
# the header styling already used by the other header cells (B1:H1)
# and filling in the per-row numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells -----------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the existing header formatting (bold font, borders, centered
# alignment - style index 1 in the original file) from H1 onto the two
# new header cells so they match the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data rows (2-38) ---------------------------------------------------
$iVals = @(7,8,6,8,7,6,5,8,8,10,4,7,9,8,9,7,1,1,6,1,1,1,10,1,6,6,6,9,6,8,7,8,5,5,6,8,3)
$jVals = @(8,8,8,8,8,8,7,8,8,10,8,8,9,9,9,8,5,3,8,5,5,5,10,5,8,9,8,10,7,8,8,9,6,8,8,9,4)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
